# Rename rule B11 from "R40" to "1".
#
# B11 currently holds the shared string "R40". The target value is the
# text "1" -- a *string*, not a number (it stays a shared-string cell,
# same style). Assigning $rng.Value = "1" directly would be re-parsed as
# user input and coerced into a numeric cell (Excel's normal "looks like a
# number -> store as number" behavior), which is not what we want here.
#
# So: write a formula that evaluates to the text "1", then convert that
# formula result to a static value via Copy/PasteSpecial (values only).
# Pasting values copies the already-typed (string) result instead of
# re-parsing "1" as user input, so the cell keeps its text type and its
# existing style/formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("B11")
$rng.Formula = '="1"'
$rng.Copy()
$rng.PasteSpecial(-4163)          # xlPasteValues
$excel.CutCopyMode = $false
